# Refresh Universalis market-price derived columns (H:N) across the Tonberry
# Profits leve tables -- scheduled-runner sync of currentAveragePrice* /
# LevePrice*/LeveProfit* figures per sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2079.875
$ws.Range("I40").Value = 2187.5
$ws.Range("J40").Value = 1972.25
$ws.Range("K40").Value = 2187.5
$ws.Range("L40").Value = 1972.25
$ws.Range("M40").Value = -2012.5
$ws.Range("N40").Value = -2322.25
$ws.Range("H41").Value = 534.125
$ws.Range("I41").Value = 209
$ws.Range("J41").Value = 609.1539
$ws.Range("K41").Value = 209
$ws.Range("L41").Value = 609.1539
$ws.Range("M41").Value = 231
$ws.Range("N41").Value = -1489.1539
$ws.Range("H62").Value = 2942.75
$ws.Range("I62").Value = 2908.4
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2908.4
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2284.4
$ws.Range("N62").Value = -4248
$ws.Range("H64").Value = 3221.4285
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 3310
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 3310
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -3806
$ws.Range("H65").Value = 2942.75
$ws.Range("I65").Value = 2908.4
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 14542
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -11422
$ws.Range("N65").Value = -21240
$ws.Range("H67").Value = 3221.4285
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 3310
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 3310
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -5026
$ws.Range("H92").Value = 2462791.5
$ws.Range("I92").Value = 3078240
$ws.Range("J92").Value = 997
$ws.Range("K92").Value = 3078240
$ws.Range("L92").Value = 997
$ws.Range("M92").Value = -3076992
$ws.Range("N92").Value = -3493
$ws.Range("H101").Value = 995.75
$ws.Range("I101").Value = 327.66666
$ws.Range("J101").Value = 3000
$ws.Range("K101").Value = 982.9999799999999
$ws.Range("L101").Value = 9000
$ws.Range("M101").Value = 639.0000200000001
$ws.Range("N101").Value = -12244
$ws.Range("H125").Value = 1074.2
$ws.Range("I125").Value = 1116.4
$ws.Range("J125").Value = 1032
$ws.Range("K125").Value = 10047.6
$ws.Range("L125").Value = 9288
$ws.Range("M125").Value = -7587.6
$ws.Range("N125").Value = -14208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2781.397
$ws.Range("I32").Value = 2089.776
$ws.Range("J32").Value = 6792.8
$ws.Range("K32").Value = 2089.776
$ws.Range("L32").Value = 6792.8
$ws.Range("M32").Value = -1802.776
$ws.Range("N32").Value = -7366.8
$ws.Range("H132").Value = 1893.4445
$ws.Range("I132").Value = 1277.1875
$ws.Range("J132").Value = 2789.818
$ws.Range("K132").Value = 3831.5625
$ws.Range("L132").Value = 8369.454000000002
$ws.Range("M132").Value = -1301.5625
$ws.Range("N132").Value = -13429.454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2253.6667
$ws.Range("I105").Value = 2214.9546
$ws.Range("J105").Value = 2424
$ws.Range("K105").Value = 2214.9546
$ws.Range("L105").Value = 2424
$ws.Range("M105").Value = -467.9546
$ws.Range("N105").Value = -5918

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1072
$ws.Range("I22").Value = 350
$ws.Range("J22").Value = 1278.2858
$ws.Range("K22").Value = 350
$ws.Range("L22").Value = 1278.2858
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = -1978.2858
$ws.Range("H31").Value = 1523.5652
$ws.Range("I31").Value = 1296.0555
$ws.Range("J31").Value = 2342.6
$ws.Range("K31").Value = 1296.0555
$ws.Range("L31").Value = 2342.6
$ws.Range("M31").Value = -1001.0555
$ws.Range("N31").Value = -2932.6
$ws.Range("H34").Value = 1523.5652
$ws.Range("I34").Value = 1296.0555
$ws.Range("J34").Value = 2342.6
$ws.Range("K34").Value = 1296.0555
$ws.Range("L34").Value = 2342.6
$ws.Range("M34").Value = -1094.0555
$ws.Range("N34").Value = -2746.6
$ws.Range("H58").Value = 967480.75
$ws.Range("I58").Value = 1318331
$ws.Range("J58").Value = 2642.5833
$ws.Range("K58").Value = 1318331
$ws.Range("L58").Value = 2642.5833
$ws.Range("M58").Value = -1318128
$ws.Range("N58").Value = -3048.5833
$ws.Range("H136").Value = 967480.75
$ws.Range("I136").Value = 1318331
$ws.Range("J136").Value = 2642.5833
$ws.Range("K136").Value = 3954993
$ws.Range("L136").Value = 7927.749899999999
$ws.Range("M136").Value = -3952443
$ws.Range("N136").Value = -13027.7499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 22681.4
$ws.Range("I55").Value = 50304
$ws.Range("J55").Value = 4266.3335
$ws.Range("K55").Value = 150912
$ws.Range("L55").Value = 12799.0005
$ws.Range("M55").Value = -150735
$ws.Range("N55").Value = -13153.0005
$ws.Range("H131").Value = 11189.045
$ws.Range("I131").Value = 491.8889
$ws.Range("J131").Value = 12848.948
$ws.Range("K131").Value = 1475.6667
$ws.Range("L131").Value = 38546.844
$ws.Range("M131").Value = 3564.3333
$ws.Range("N131").Value = -48626.844
$ws.Range("H133").Value = 1998
$ws.Range("I133").Value = 1998
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 5994
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -934
$ws.Range("N133").ClearContents()
$ws.Range("H137").Value = 4387
$ws.Range("I137").Value = 2572.5
$ws.Range("J137").Value = 5294.25
$ws.Range("K137").Value = 7717.5
$ws.Range("L137").Value = 15882.75
$ws.Range("M137").Value = -2617.5
$ws.Range("N137").Value = -26082.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3558.111
$ws.Range("I102").Value = 3789
$ws.Range("J102").Value = 2750
$ws.Range("K102").Value = 3789
$ws.Range("L102").Value = 2750
$ws.Range("M102").Value = -2167
$ws.Range("N102").Value = -5994
$ws.Range("H132").Value = 876101.9399999999
$ws.Range("I132").Value = 1375275.5
$ws.Range("J132").Value = 2548.1875
$ws.Range("K132").Value = 4125826.5
$ws.Range("L132").Value = 7644.5625
$ws.Range("M132").Value = -4123296.5
$ws.Range("N132").Value = -12704.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 264.7
$ws.Range("I55").Value = 183.375
$ws.Range("J55").Value = 590
$ws.Range("K55").Value = 183.375
$ws.Range("L55").Value = 590
$ws.Range("M55").Value = -10.375
$ws.Range("N55").Value = -936
$ws.Range("H122").Value = 7738.8
$ws.Range("I122").Value = 1704
$ws.Range("J122").Value = 9247.5
$ws.Range("K122").Value = 5112
$ws.Range("L122").Value = 27742.5
$ws.Range("M122").Value = -2662
$ws.Range("N122").Value = -32642.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1733.7273
$ws.Range("I81").Value = 1783.875
$ws.Range("J81").Value = 1600
$ws.Range("K81").Value = 3567.75
$ws.Range("L81").Value = 3200
$ws.Range("M81").Value = -2506.75
$ws.Range("N81").Value = -5322
$ws.Range("H84").Value = 1733.7273
$ws.Range("I84").Value = 1783.875
$ws.Range("J84").Value = 1600
$ws.Range("K84").Value = 17838.75
$ws.Range("L84").Value = 16000
$ws.Range("M84").Value = -12534.75
$ws.Range("N84").Value = -26608
$ws.Range("H108").Value = 63999
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 63999
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 63999
$ws.Range("N108").Value = -71679
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H122").Value = 66015.086
$ws.Range("I122").Value = 66015.086
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 198045.258
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -195595.258
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 10785.857
$ws.Range("I126").Value = 11000.182
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 33000.546
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -30530.546
$ws.Range("N126").Value = -34940
$ws.Range("H132").Value = 1332.7742
$ws.Range("I132").Value = 972.52
$ws.Range("J132").Value = 2833.8333
$ws.Range("K132").Value = 2917.56
$ws.Range("L132").Value = 8501.499899999999
$ws.Range("M132").Value = -387.5599999999999
$ws.Range("N132").Value = -13561.4999
$ws.Range("H136").Value = 18521192
$ws.Range("I136").Value = 25255188
$ws.Range("J136").Value = 2698.75
$ws.Range("K136").Value = 75765564
$ws.Range("L136").Value = 8096.25
$ws.Range("M136").Value = -75763014
$ws.Range("N136").Value = -13196.25

